{"js": "const replacements = [\n  [\"15\u00d724=360\", \"44\u00d737=1628\"],\n  [\"31\u00d789=2759\", \"91\u00d753=4823\"],\n  [\"85\u00d711=935\", \"72\u00d786=6192\"],\n  [\"55\u00d728=1540\", \"54\u00d791=4914\"],\n  [\"63\u00d735=2205\", \"31\u00d730=930\"],\n  [\"21\u00d713=273\", \"53\u00d778=4134\"],\n  [\"57\u00d799=5643\", \"73\u00d789=6497\"],\n  [\"38\u00d732=1216\", \"22\u00d734=748\"],\n  [\"92\u00d784=7728\", \"31\u00d764=1984\"],\n  [\"54\u00d773=3942\", \"48\u00d787=4176\"],\n  [\"52\u00d730=1560\", \"78\u00d741=3198\"],\n  [\"65\u00d719=1235\", \"31\u00d714=434\"],\n  [\"17\u00d759=1003\", \"64\u00d773=4672\"],\n  [\"40\u00d718=720\", \"48\u00d795=4560\"],\n  [\"21\u00d770=1470\", \"86\u00d714=1204\"],\n  [\"42\u00d777=3234\", \"20\u00d785=1700\"],\n  [\"76\u00d785=6460\", \"56\u00d795=5320\"],\n  [\"18\u00d795=1710\", \"65\u00d740=2600\"],\n  [\"31\u00d758=1798\", \"17\u00d762=1054\"],\n  [\"41\u00d793=3813\", \"95\u00d756=5320\"],\n  [\"79\u00d715=1185\", \"42\u00d727=1134\"],\n  [\"58\u00d727=1566\", \"61\u00d750=3050\"],\n  [\"70\u00d733=2310\", \"76\u00d751=3876\"],\n  [\"29\u00d773=2117\", \"51\u00d778=3978\"],\n  [\"88\u00d731=2728\", \"14\u00d712=168\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"15\u00d724=360\", \"44\u00d737=1628\"),\n  @(\"31\u00d789=2759\", \"91\u00d753=4823\"),\n  @(\"85\u00d711=935\", \"72\u00d786=6192\"),\n  @(\"55\u00d728=1540\", \"54\u00d791=4914\"),\n  @(\"63\u00d735=2205\", \"31\u00d730=930\"),\n  @(\"21\u00d713=273\", \"53\u00d778=4134\"),\n  @(\"57\u00d799=5643\", \"73\u00d789=6497\"),\n  @(\"38\u00d732=1216\", \"22\u00d734=748\"),\n  @(\"92\u00d784=7728\", \"31\u00d764=1984\"),\n  @(\"54\u00d773=3942\", \"48\u00d787=4176\"),\n  @(\"52\u00d730=1560\", \"78\u00d741=3198\"),\n  @(\"65\u00d719=1235\", \"31\u00d714=434\"),\n  @(\"17\u00d759=1003\", \"64\u00d773=4672\"),\n  @(\"40\u00d718=720\", \"48\u00d795=4560\"),\n  @(\"21\u00d770=1470\", \"86\u00d714=1204\"),\n  @(\"42\u00d777=3234\", \"20\u00d785=1700\"),\n  @(\"76\u00d785=6460\", \"56\u00d795=5320\"),\n  @(\"18\u00d795=1710\", \"65\u00d740=2600\"),\n  @(\"31\u00d758=1798\", \"17\u00d762=1054\"),\n  @(\"41\u00d793=3813\", \"95\u00d756=5320\"),\n  @(\"79\u00d715=1185\", \"42\u00d727=1134\"),\n  @(\"58\u00d727=1566\", \"61\u00d750=3050\"),\n  @(\"70\u00d733=2310\", \"76\u00d751=3876\"),\n  @(\"29\u00d773=2117\", \"51\u00d778=3978\"),\n  @(\"88\u00d731=2728\", \"14\u00d712=168\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
